$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 9573.909
$ws.Range("I12").Value = 267.8
$ws.Range("K12").Value = 267.8
$ws.Range("M12").Value = -97.80000000000001
$ws.Range("H41").Value = 1043.5714
$ws.Range("I41").Value = 311
$ws.Range("J41").Value = 2020.3334
$ws.Range("K41").Value = 311
$ws.Range("L41").Value = 2020.3334
$ws.Range("M41").Value = 129
$ws.Range("N41").Value = -2900.3334
$ws.Range("H43").Value = 6937.125
$ws.Range("J43").Value = 6785.4287
$ws.Range("L43").Value = 6785.4287
$ws.Range("N43").Value = -6923.4287
$ws.Range("H113").Value = 9635.875
$ws.Range("J113").Value = 6324
$ws.Range("L113").Value = 6324
$ws.Range("N113").Value = -12832
$ws.Range("H116").Value = 1393781
$ws.Range("I116").Value = 2224652
$ws.Range("J116").Value = 8996
$ws.Range("K116").Value = 2224652
$ws.Range("L116").Value = 8996
$ws.Range("M116").Value = -2221210
$ws.Range("N116").Value = -15880
$ws.Range("H125").Value = 4775.9033
$ws.Range("J125").Value = 4886.931
$ws.Range("L125").Value = 43982.37899999999
$ws.Range("N125").Value = -48902.37899999999
$ws.Range("H137").Value = 271441.9
$ws.Range("I137").Value = 377226.22
$ws.Range("J137").Value = 2172.7273
$ws.Range("K137").Value = 1131678.66
$ws.Range("L137").Value = 6518.1819
$ws.Range("M137").Value = -1129128.66
$ws.Range("N137").Value = -11618.1819

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6079.6665
$ws.Range("I45").Value = 5210.778
$ws.Range("J45").Value = 7383
$ws.Range("K45").Value = 5210.778
$ws.Range("L45").Value = 7383
$ws.Range("M45").Value = -4833.778
$ws.Range("N45").Value = -8137
$ws.Range("H46").Value = 8033.75
$ws.Range("J46").Value = 8033.75
$ws.Range("L46").Value = 8033.75
$ws.Range("N46").Value = -8671.75
$ws.Range("H61").Value = 3700.4
$ws.Range("I61").Value = 2981.9656
$ws.Range("J61").Value = 5594.4546
$ws.Range("K61").Value = 2981.9656
$ws.Range("L61").Value = 5594.4546
$ws.Range("M61").Value = -2769.9656
$ws.Range("N61").Value = -6018.4546
$ws.Range("H74").Value = 2202.2363
$ws.Range("I74").Value = 1237.7059
$ws.Range("K74").Value = 1237.7059
$ws.Range("M74").Value = -363.7058999999999
$ws.Range("H77").Value = 2202.2363
$ws.Range("I77").Value = 1237.7059
$ws.Range("K77").Value = 6188.5295
$ws.Range("M77").Value = -1820.5295
$ws.Range("H97").Value = 10127.857
$ws.Range("I97").Value = 10707.923
$ws.Range("K97").Value = 10707.923
$ws.Range("M97").Value = -10211.923
$ws.Range("H102").Value = 4972.2666
$ws.Range("I102").Value = 3603.238
$ws.Range("K102").Value = 3603.238
$ws.Range("M102").Value = -1981.238
$ws.Range("H110").Value = 7327.0527
$ws.Range("I110").Value = 9247.308000000001
$ws.Range("K110").Value = 9247.308000000001
$ws.Range("M110").Value = -7202.308000000001
$ws.Range("H132").Value = 3871.0435
$ws.Range("J132").Value = 5293.875
$ws.Range("L132").Value = 15881.625
$ws.Range("N132").Value = -20941.625
$ws.Range("H133").Value = 57750
$ws.Range("J133").Value = 57750
$ws.Range("L133").Value = 57750
$ws.Range("N133").Value = -62810
$ws.Range("H136").Value = 3700.4
$ws.Range("I136").Value = 2981.9656
$ws.Range("J136").Value = 5594.4546
$ws.Range("K136").Value = 8945.8968
$ws.Range("L136").Value = 16783.3638
$ws.Range("M136").Value = -6395.8968
$ws.Range("N136").Value = -21883.3638

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 32092.834
$ws.Range("I99").Value = 34815.184
$ws.Range("J99").Value = 2147
$ws.Range("K99").Value = 34815.184
$ws.Range("L99").Value = 2147
$ws.Range("M99").Value = -33317.184
$ws.Range("N99").Value = -5143
$ws.Range("H134").Value = 2720.8923
$ws.Range("I134").Value = 2019
$ws.Range("J134").Value = 4553.6113
$ws.Range("K134").Value = 6057
$ws.Range("L134").Value = 13660.8339
$ws.Range("M134").Value = -3522
$ws.Range("N134").Value = -18730.8339

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 35991356
$ws.Range("I6").Value = 35991356
$ws.Range("K6").Value = 35991356
$ws.Range("M6").Value = -35991243
$ws.Range("H22").Value = 999.75
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("M22").Value = -650
$ws.Range("H58").Value = 3186.2666
$ws.Range("I58").Value = 2039.7222
$ws.Range("K58").Value = 2039.7222
$ws.Range("M58").Value = -1836.7222
$ws.Range("H134").Value = 3680932
$ws.Range("I134").Value = 3680932
$ws.Range("K134").Value = 11042796
$ws.Range("M134").Value = -11040261
$ws.Range("H136").Value = 3186.2666
$ws.Range("I136").Value = 2039.7222
$ws.Range("K136").Value = 6119.1666
$ws.Range("M136").Value = -3569.1666

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 799.6667
$ws.Range("I3").Value = 700
$ws.Range("K3").Value = 2100
$ws.Range("M3").Value = -1988
$ws.Range("H4").Value = 15113147
$ws.Range("I4").Value = 2516996
$ws.Range("K4").Value = 7550988
$ws.Range("M4").Value = -7550876
$ws.Range("H133").Value = 4410
$ws.Range("I133").Value = 3538.182
$ws.Range("K133").Value = 10614.546
$ws.Range("M133").Value = -5554.545999999998

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3199
$ws.Range("I80").Value = 2832.6667
$ws.Range("J80").Value = 3473.75
$ws.Range("K80").Value = 2832.6667
$ws.Range("L80").Value = 3473.75
$ws.Range("M80").Value = -1834.6667
$ws.Range("N80").Value = -5469.75
$ws.Range("H83").Value = 3199
$ws.Range("I83").Value = 2832.6667
$ws.Range("J83").Value = 3473.75
$ws.Range("K83").Value = 14163.3335
$ws.Range("L83").Value = 17368.75
$ws.Range("M83").Value = -9171.333500000001
$ws.Range("N83").Value = -27352.75
$ws.Range("H99").Value = 17984.75
$ws.Range("J99").Value = 27499.5
$ws.Range("L99").Value = 27499.5
$ws.Range("N99").Value = -31991.5
$ws.Range("H123").Value = 18021.695
$ws.Range("J123").Value = 18021.695
$ws.Range("L123").Value = 18021.695
$ws.Range("N123").Value = -22921.695
$ws.Range("H132").Value = 2877.4827
$ws.Range("I132").Value = 2877.4827
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8632.4481
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6102.4481
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 100019.7
$ws.Range("J136").Value = 100019.7
$ws.Range("L136").Value = 300059.1
$ws.Range("N136").Value = -305159.1

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6266.6875
$ws.Range("I136").Value = 2398.1428
$ws.Range("J136").Value = 9275.556
$ws.Range("K136").Value = 7194.428400000001
$ws.Range("L136").Value = 27826.668
$ws.Range("M136").Value = -4644.428400000001
$ws.Range("N136").Value = -32926.66800000001

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 260583.5
$ws.Range("I62").Value = 284000.2
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 284000.2
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -283376.2
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 260583.5
$ws.Range("I65").Value = 284000.2
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 1420001
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -1416881
$ws.Range("N65").Value = -21240
$ws.Range("H107").Value = 41999.625
$ws.Range("I107").Value = 4199.6
$ws.Range("J107").Value = 104999.664
$ws.Range("K107").Value = 12598.8
$ws.Range("L107").Value = 314998.992
$ws.Range("M107").Value = -10678.8
$ws.Range("N107").Value = -318838.992
$ws.Range("H126").Value = 25627.316
$ws.Range("I126").Value = 33804.92
$ws.Range("K126").Value = 101414.76
$ws.Range("M126").Value = -98944.75999999999
$ws.Range("H136").Value = 2763.182
$ws.Range("J136").Value = 3510.75
$ws.Range("L136").Value = 10532.25
$ws.Range("N136").Value = -15632.25
